$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 119 (shifts old rows 119:242 down to 120:243)
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with its data
# (a new weekly price record for "Pepino dulce" / Cultivar IV Región / Primera)
$ws.Cells.Item(119, 1).Value = 10
$ws.Cells.Item(119, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(119, 3).Value = "La Araucanía"
$ws.Cells.Item(119, 4).Value = 44740
$ws.Cells.Item(119, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(119, 5).Value = 9
$ws.Cells.Item(119, 6).Value = 100112043
$ws.Cells.Item(119, 7).Value = "Pepino dulce"
$ws.Cells.Item(119, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 200
$ws.Cells.Item(119, 11).Value = 19000
$ws.Cells.Item(119, 12).Value = 19000
$ws.Cells.Item(119, 13).Value = 19000
$ws.Cells.Item(119, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(119, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(119, 16).Value = 1056
$ws.Cells.Item(119, 17).Value = 18
$ws.Cells.Item(119, 18).Value = "Hortaliza"
